$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (price) cells: force text so dotted numbers are not reinterpreted as numbers ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.536.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.766.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.763.51"
$ws.Range("D7").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.40"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.463"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000251"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.399.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.767.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.538.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "497.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.732"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000154"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.06"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.916.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.703.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "441.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.815.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "140.86"
$ws.Range("D50").Style = "Normal"

# --- Column B/C/E cells: plain text assignment ---
$ws.Range("E2").Value = "  +3.19%  "
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("E6").Value = "  +2.48%  "
$ws.Range("E7").Value = "  +1.54%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +2.21%  "
$ws.Range("E10").Value = "  +4.25%  "
$ws.Range("E11").Value = "  +3.80%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("E13").Value = "  +2.60%  "
$ws.Range("E14").Value = "  +4.40%  "
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("E17").Value = "  +3.03%  "
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("E21").Value = "  +15.97%  "
$ws.Range("E22").Value = "  +2.09%  "
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("E24").Value = "  +11.75%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  +2.51%  "
$ws.Range("E27").Value = "  +2.25%  "
$ws.Range("E28").Value = "  +2.34%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("E30").Value = "  +2.97%  "
$ws.Range("E31").Value = "  +6.83%  "
$ws.Range("E32").Value = "  +5.92%  "
$ws.Range("E33").Value = "  +1.89%  "
$ws.Range("E34").Value = "  +1.85%  "
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("E38").Value = "  +2.05%  "
$ws.Range("E39").Value = "  +3.23%  "
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("E41").Value = "  +1.73%  "
$ws.Range("E42").Value = "  +10.57%  "
$ws.Range("E43").Value = "  +1.99%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("E44").Value = "  +3.31%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("E49").Value = "  +2.22%  "
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("E51").Value = "  +2.76%  "
